$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.873.88"
$ws.Range("E2").Value = "  +6.43%  "

$ws.Range("D3").Value = "3.880.02"
$ws.Range("E3").Value = "  +13.23%  "

$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "'425.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.36%  "

$ws.Range("D6").Value = "'130.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.93%  "

$ws.Range("D7").Value = "3.875.80"
$ws.Range("E7").Value = "  +13.40%  "

$ws.Range("D8").Value = "'0.610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.50%  "

$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").Value = "'0.722"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.30%  "

$ws.Range("E11").Value = "  +13.77%  "

$ws.Range("D12").Value = "'0.0000343"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +18.51%  "

$ws.Range("D13").Value = "'40.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.26%  "

$ws.Range("D14").Value = "4.498.56"
$ws.Range("E14").Value = "  +13.67%  "

$ws.Range("E15").Value = "  +12.96%  "

$ws.Range("D16").Value = "'15.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +31.68%  "

$ws.Range("D17").Value = "3.931.28"
$ws.Range("E17").Value = "  +14.93%  "

$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("D19").Value = "'19.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.92%  "

$ws.Range("D20").Value = "67.164.17"
$ws.Range("E20").Value = "  +6.91%  "

$ws.Range("E21").Value = "  +7.95%  "

$ws.Range("D22").Value = "'411.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.85%  "

$ws.Range("D23").Value = "'14.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.75%  "

$ws.Range("D24").Value = "'84.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.99%  "

$ws.Range("E25").Value = "  +9.61%  "

$ws.Range("D26").Value = "'37.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.92%  "

$ws.Range("D27").Value = "'3.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.19%  "

$ws.Range("E28").Value = "  +15.29%  "

$ws.Range("E29").Value = "  +3.70%  "

$ws.Range("D30").Value = "'8.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +35.79%  "

$ws.Range("D31").Value = "'729.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.69%  "

$ws.Range("D32").Value = "'13.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.67%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "'2.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.57%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.121"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.21%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "'38.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.91%  "

$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("D38").Value = "'55.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.04%  "

$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'5.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +34.44%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0756"
$ws.Range("E40").Value = "  +29.40%  "

$ws.Range("D41").Value = "'0.0458"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.61%  "

$ws.Range("D42").Value = "'2.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.72%  "

$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.00%  "

$ws.Range("E45").Value = "  +4.00%  "

$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").Value = "'3.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.51%  "

$ws.Range("E47").Value = "  +15.53%  "

$ws.Range("D48").Value = "'142.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.31%  "

$ws.Range("E49").Value = "  +7.46%  "

$ws.Range("E50").Value = "  +8.21%  "

$ws.Range("E51").Value = "  +5.93%  "
